# Apply the changes described by the diff:
# 1. Rename header "Requested quantity" -> "Weekly_PO_Qty" on "Weekly Quantity" sheet (B1)
# 2. Rename header "Requested quantity" -> "Monthly_PO_Qty" on "Monthly Trend" sheet (B1)
# 3. Add a new "PO Forecast" worksheet (after the existing sheets) with forecast data

$wb = $excel.ActiveWorkbook

# --- 1 & 2: update header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$poData = @(
    @(44941.99999999999, 275, 93.41071888335037, 460.2326808143995),
    @(44983.99999999999, 217, 32.85665079355153, 399.3156969030751),
    @(45011.99999999999, 178, -2.119385226610473, 367.9231862279708),
    @(45018.99999999999, 168, -11.31047828924579, 351.4906547800385),
    @(45025.99999999999, 159, -21.25680093217206, 348.9623699958609),
    @(45046.99999999999, 130, -42.73166069028922, 313.2733032493877),
    @(45053.99999999999, 120, -68.79713705559487, 298.2872876101627),
    @(45060.99999999999, 110, -68.63753916941161, 299.4116694605044),
    @(45067.99999999999, 100, -68.60878497798195, 279.0918704509562),
    @(45074.99999999999, 91, -101.0922497599437, 290.4085391871518),
    @(45081.99999999999, 81, -97.21495122205083, 273.9912263312399),
    @(45088.99999999999, 71, -116.0204751610123, 255.6745412639911),
    @(45095.99999999999, 62, -126.2679975513668, 249.5218262655057),
    @(45102.99999999999, 52, -144.4009975264809, 236.6094037534338),
    @(45109.99999999999, 42, -127.9359430133345, 244.565901380842),
    @(45116.99999999999, 33, -153.1789078061049, 215.4274363773978),
    @(45123.99999999999, 23, -150.3585970363426, 213.0587825959998)
)

$r = 2
foreach ($row in $poData) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

Write-Output "PO Forecast sheet added; headers updated."
